$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 4
$ws.Range("C4").Value = 0.4302910291029102
$ws.Range("D4").Value = 0.9818301830183017
$ws.Range("E4").Value = -0.9999879987998799
$ws.Range("F4").Value = 0.9783978397839782
$ws.Range("H4").Value = 0.7083228322832282
$ws.Range("I4").Value = -0.01120912091209121
$ws.Range("J4").Value = -0.1741974197419742
$ws.Range("K4").Value = 0.1503150315031503
$ws.Range("L4").Value = -0.0722112211221122
$ws.Range("M4").Value = -0.08854485448544852
$ws.Range("N4").Value = 0.2058205820582058
$ws.Range("P4").Value = -0.06069006900690068

# Row 5
$ws.Range("C5").Value = -0.04244824482448245
$ws.Range("D5").Value = -0.1724452445244524
$ws.Range("E5").Value = 0.1723372337233723
$ws.Range("F5").Value = -0.1661326132613261
$ws.Range("H5").Value = -0.1723852385238524
$ws.Range("I5").Value = -0.1615121512151215
$ws.Range("J5").Value = 0.1104950495049505
$ws.Range("K5").Value = 0.1734893489348935
$ws.Range("L5").Value = -0.02035403540354035
$ws.Range("M5").Value = -0.07474347434743474
$ws.Range("N5").Value = -0.006396639663966396
$ws.Range("P5").Value = 0.03697569756975697

# Row 6
$ws.Range("C6").Value = 0.1437743774377438
$ws.Range("D6").Value = 0.09165316531653164
$ws.Range("E6").Value = -0.1153675367536754
$ws.Range("F6").Value = 0.1654725472547255
$ws.Range("H6").Value = 0.09804980498049803
$ws.Range("I6").Value = -0.05526552655265525
$ws.Range("J6").Value = 0.1964476447644764
$ws.Range("K6").Value = -0.04043204320432043
$ws.Range("L6").Value = -0.1233483348334833
$ws.Range("M6").Value = -0.05747374737473747
$ws.Range("N6").Value = -0.06136213621362136
$ws.Range("P6").Value = -0.03842784278427842

# Row 7
$ws.Range("C7").Value = 0.1685208520852085
$ws.Range("D7").Value = 0.1241284128412841
$ws.Range("E7").Value = 0.04049204920492049
$ws.Range("F7").Value = -0.154035403540354
$ws.Range("H7").Value = -0.08687668766876687
$ws.Range("I7").Value = 0.004116411641164116
$ws.Range("J7").Value = 0.1304530453045304
$ws.Range("K7").Value = 0.2242184218421842
$ws.Range("L7").Value = 0.9638283828382838
$ws.Range("M7").Value = -0.01608160816081608
$ws.Range("N7").Value = -0.644044404440444
$ws.Range("P7").Value = -0.1844824482448245

# Row 8
$ws.Range("C8").Value = 0.0814161416141614
$ws.Range("D8").Value = -0.05588958895889588
$ws.Range("E8").Value = 0.06711071107110711
$ws.Range("F8").Value = 0.01669366936693669
$ws.Range("H8").Value = 0.638163816381638
$ws.Range("I8").Value = 0.7969636963696368
$ws.Range("J8").Value = -0.07981998199819981
$ws.Range("K8").Value = 0.02159015901590159
$ws.Range("L8").Value = 0.004872487248724872
$ws.Range("M8").Value = 0.4324872487248725
$ws.Range("N8").Value = 0.117995799579958
$ws.Range("P8").Value = 0.2530573057305731

# Row 9
$ws.Range("C9").Value = 0.4402880288028803
$ws.Range("D9").Value = 0.02558655865586558
$ws.Range("E9").Value = -0.04106810681068106
$ws.Range("F9").Value = 0.03251125112511251
$ws.Range("H9").Value = 0.02144614461446144
$ws.Range("I9").Value = 0.1024542454245424
$ws.Range("J9").Value = -0.09371737173717372
$ws.Range("K9").Value = -0.6194179417941794
$ws.Range("L9").Value = -0.03203120312031203
$ws.Range("M9").Value = -0.06874287428742873
$ws.Range("N9").Value = 0.08735673567356733
$ws.Range("P9").Value = -0.05572157215721572

# Row 10
$ws.Range("C10").Value = 0.484032403240324
$ws.Range("D10").Value = -0.05312931293129312
$ws.Range("E10").Value = 0.06724272427242724
$ws.Range("F10").Value = -0.0594059405940594
$ws.Range("H10").Value = -0.00936093609360936
$ws.Range("I10").Value = -0.05666966696669666
$ws.Range("J10").Value = -0.03473147314731473
$ws.Range("K10").Value = 0.5953195319531953
$ws.Range("L10").Value = 0.07278727872787277
$ws.Range("M10").Value = 0.04853285328532853
$ws.Range("N10").Value = 0.01035703570357036
$ws.Range("P10").Value = -0.01855385538553855

# Row 11
$ws.Range("C11").Value = -0.1516471647164716
$ws.Range("D11").Value = 0.07734773477347733
$ws.Range("E11").Value = -0.09911791179117913
$ws.Range("F11").Value = 0.1065826582658266
$ws.Range("H11").Value = -0.05977797779777976
$ws.Range("I11").Value = -0.1155475547554755
$ws.Range("J11").Value = 0.04516051605160516
$ws.Range("K11").Value = -0.07468346834683467
$ws.Range("L11").Value = -0.1444944494449445
$ws.Range("M11").Value = -0.01642964296429643
$ws.Range("N11").Value = 0.0372997299729973
$ws.Range("P11").Value = 0.000204020402040204

# Row 12
$ws.Range("C12").Value = -0.03000300030003
$ws.Range("D12").Value = -0.1913831383138313
$ws.Range("E12").Value = 0.1782058205820582
$ws.Range("F12").Value = -0.1818541854185418
$ws.Range("H12").Value = -0.1975397539753975
$ws.Range("I12").Value = -0.1073267326732673
$ws.Range("J12").Value = 0.04489648964896489
$ws.Range("K12").Value = 0.04469246924692469
$ws.Range("L12").Value = 0.04136813681368137
$ws.Range("M12").Value = 0.07413141314131412
$ws.Range("N12").Value = -0.0101050105010501
$ws.Range("P12").Value = -0.08259225922592257

# Row 13
$ws.Range("C13").Value = 0.2233303330333033
$ws.Range("D13").Value = 0.07692769276927691
$ws.Range("E13").Value = -0.08200420042004199
$ws.Range("F13").Value = 0.07575157515751575
$ws.Range("H13").Value = 0.03144314431443144
$ws.Range("I13").Value = 0.564128412841284
$ws.Range("J13").Value = 0.05083708370837083
$ws.Range("K13").Value = -0.1886108610861086
$ws.Range("L13").Value = -0.004908490849084907
$ws.Range("M13").Value = 0.1396939693969397
$ws.Range("N13").Value = -0.005880588058805879
$ws.Range("P13").Value = -0.114047404740474

# Row 14
$ws.Range("C14").Value = -0.3244404440444044
$ws.Range("D14").Value = -0.01903390339033903
$ws.Range("E14").Value = -0.003084308430843084
$ws.Range("F14").Value = -0.003228322832283228
$ws.Range("H14").Value = 0.04126012601260126
$ws.Range("I14").Value = -0.05876987698769877
$ws.Range("J14").Value = -0.1136513651365136
$ws.Range("K14").Value = -0.1269606960696069
$ws.Range("L14").Value = -0.1572397239723972
$ws.Range("M14").Value = 0.0909210921092109
$ws.Range("N14").Value = 0.09971797179717971
$ws.Range("P14").Value = -0.05868586858685867
